# The document has a collapsed "_GoBack" bookmark sitting right after the
# text of the "pb,市净率" paragraph. This edit relocates that bookmark so it
# instead sits right after the text of the "gpr,毛利率(%)" paragraph
# (six paragraphs further down the body) - i.e. immediately before the
# closing paragraph mark of that paragraph, exactly as it previously sat
# relative to the "pb" paragraph.

$d = $word.ActiveDocument

# 1. Remove the bookmark from wherever it currently is.
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

# 2. Locate the "gpr,毛利率(%)" paragraph (the new home for the bookmark)
#    by its text rather than a hard-coded index, so the script keeps working
#    even if paragraphs shift around.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "gpr*") {
        $targetPara = $p
    }
}

$target = $targetPara.Range
$pos = $target.End - 1

# 3. Placing a bookmark straight onto a zero-length (collapsed) Range at a
#    paragraph-end boundary isn't reliable, so instead: insert a throwaway
#    placeholder character there, wrap a (non-collapsed) bookmark around it,
#    then delete the placeholder. Word automatically keeps the bookmark in
#    place (collapsed) once its contents are deleted - the same mechanism
#    that produces the naturally-collapsed "_GoBack" bookmark.
$d.Range($pos, $pos).InsertAfter("@")
$wrap = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $wrap)
$d.Range($pos, $pos + 1).Text = ""
